# Update the "Overview" worksheet of the income-statement workbook with the
# latest reported figures (rows 11-27, columns D-H), and clean up a couple of
# rows that previously held literal "-" placeholder text but should now be
# numeric zeros (rows 15 and 23, plus the stray F16 cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Row 11 - فروش (Sales)
$ws.Range("D11").Value = 4511858
$ws.Range("E11").Value = 9200140
$ws.Range("F11").Value = 13339160
$ws.Range("G11").Value = 18651716
$ws.Range("H11").Value = 32343396

# Row 12 - بهای تمام شده کالای فروش رفته (COGS)
$ws.Range("D12").Value = -3931914
$ws.Range("E12").Value = -8390293
$ws.Range("F12").Value = -12069436
$ws.Range("G12").Value = -16853415
$ws.Range("H12").Value = -28732239

# Row 13 - سود (زیان) ناخالص (Gross profit)
$ws.Range("D13").Value = 579944
$ws.Range("E13").Value = 809847
$ws.Range("F13").Value = 1269724
$ws.Range("G13").Value = 1798301
$ws.Range("H13").Value = 3611157

# Row 14 - هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
$ws.Range("D14").Value = -266009
$ws.Range("E14").Value = -334641
$ws.Range("F14").Value = -499483
$ws.Range("G14").Value = -779401
$ws.Range("H14").Value = -1235406

# Row 15 - هزینه کاهش ارزش دریافتنی‌ها (هزینه استثنایی)
# Previously literal "-" text in every column; now numeric zeros.
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# Row 16 - خالص سایر درامدها (هزینه ها) ی عملیاتی
$ws.Range("D16").Value = 13853
$ws.Range("E16").Value = 11709
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 18152
$ws.Range("H16").Value = 1353

# Row 17 - سود (زیان) عملیاتی (Operating profit)
$ws.Range("D17").Value = 327788
$ws.Range("E17").Value = 486915
$ws.Range("F17").Value = 770241
$ws.Range("G17").Value = 1037052
$ws.Range("H17").Value = 2377104

# Row 18 - هزینه های مالی (Financial expenses)
$ws.Range("D18").Value = -18521
$ws.Range("E18").Value = -38924
$ws.Range("F18").Value = -83370
$ws.Range("G18").Value = -47582
$ws.Range("H18").Value = -43924

# Row 19 - خالص سایر درامدها و هزینه های غیرعملیاتی
$ws.Range("D19").Value = 24465
$ws.Range("E19").Value = 40757
$ws.Range("F19").Value = 42773
$ws.Range("G19").Value = 64368
$ws.Range("H19").Value = 153792

# Row 20 - سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Range("D20").Value = 333732
$ws.Range("E20").Value = 488748
$ws.Range("F20").Value = 729644
$ws.Range("G20").Value = 1053838
$ws.Range("H20").Value = 2486972

# Row 21 - مالیات (Tax)
$ws.Range("D21").Value = -58653
$ws.Range("E21").Value = -82637
$ws.Range("F21").Value = -168799
$ws.Range("G21").Value = -192768
$ws.Range("H21").Value = -400310

# Row 22 - سود (زیان) خالص عملیات در حال تداوم
$ws.Range("D22").Value = 275079
$ws.Range("E22").Value = 406111
$ws.Range("F22").Value = 560845
$ws.Range("G22").Value = 861070
$ws.Range("H22").Value = 2086662

# Row 23 - سود (زیان) عملیات متوقف شده پس از اثر مالیاتی
# Previously literal "-" text in every column; now numeric zeros.
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# Row 24 - سود (زیان) خالص (Net profit)
$ws.Range("D24").Value = 275079
$ws.Range("E24").Value = 406111
$ws.Range("F24").Value = 560845
$ws.Range("G24").Value = 861070
$ws.Range("H24").Value = 2086662

# Row 25 - سود هر سهم پس از کسر مالیات (EPS)
$ws.Range("D25").Value = 1214
$ws.Range("E25").Value = 486
$ws.Range("F25").Value = 671
$ws.Range("G25").Value = 1030
$ws.Range("H25").Value = 1805

# Row 26 - سرمایه (Capital)
$ws.Range("D26").Value = 226509
$ws.Range("E26").Value = 835820
$ws.Range("F26").Value = 835820
$ws.Range("G26").Value = 835821
$ws.Range("H26").Value = 1156190

# Row 27 - سود هر سهم بر اساس آخرین سرمایه
$ws.Range("D27").Value = 238
$ws.Range("E27").Value = 351
$ws.Range("F27").Value = 485
$ws.Range("G27").Value = 745
$ws.Range("H27").Value = 1805
